$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C22").Value = 7320
$ws.Range("C23:C27").Value = 7312
$ws.Range("C28:C179").Value = 7310
